$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($ws, $addr, $val) {
    # Force text so Excel does not reinterpret values like "0,038" as a
    # thousands-formatted number (comma treated as group separator)
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

# Remove the last data row (id 24, Won Sul-Coreano) - table shrinks from 24 to 23 entries
$ws.Rows.Item(25).Delete()

# Update remaining data rows (2-24) with refreshed currency quotes and timestamps
$ws.Range("B2").Value = "Dólar"
$ws.Range("C2").Value = "$"
$ws.Range("D2").Value = "USD"
$ws.Range("E2").Value = "5,50"
$ws.Range("F2").Value = "20 de jun., 14:50 UTC ·"

$ws.Range("B3").Value = "Euro"
$ws.Range("C3").Value = "€"
$ws.Range("D3").Value = "EUR"
$ws.Range("E3").Value = "6,32"
$ws.Range("F3").Value = "20 de jun., 14:48 UTC ·"

$ws.Range("B4").Value = "Libra Esterlina"
$ws.Range("C4").Value = "£"
$ws.Range("D4").Value = "GBP"
$ws.Range("E4").Value = "7,40"
$ws.Range("F4").Value = "20 de jun., 14:49 UTC ·"

$ws.Range("B5").Value = "Iene"
$ws.Range("C5").Value = "¥"
$ws.Range("D5").Value = "JPY"
Set-TextCell $ws "E5" "0,038"
$ws.Range("F5").Value = "20 de jun., 14:49 UTC ·"

$ws.Range("B6").Value = "Franco Suíço"
$ws.Range("C6").Value = "CHF"
$ws.Range("D6").Value = "CHF"
$ws.Range("E6").Value = "6,72"
$ws.Range("F6").Value = "20 de jun., 14:48 UTC ·"

$ws.Range("B7").Value = "Dólar Australiano"
$ws.Range("C7").Value = "A$"
$ws.Range("D7").Value = "AUD"
$ws.Range("E7").Value = "3,55"
$ws.Range("F7").Value = "20 de jun., 14:48 UTC ·"

$ws.Range("B8").Value = "Peso Mexicano"
$ws.Range("C8").Value = "$"
$ws.Range("D8").Value = "MXN"
$ws.Range("E8").Value = "0,29"
$ws.Range("F8").Value = "20 de jun., 14:48 UTC ·"

$ws.Range("B9").Value = "Dólar Canadiano"
$ws.Range("C9").Value = "$"
$ws.Range("D9").Value = "CAD"
$ws.Range("E9").Value = "4,00"
$ws.Range("F9").Value = "20 de jun., 14:48 UTC ·"

$ws.Range("B10").Value = "Dólar de Hong Kong"
$ws.Range("C10").Value = "HK$"
$ws.Range("D10").Value = "HKD"
$ws.Range("E10").Value = "0,70"
$ws.Range("F10").Value = "20 de jun., 14:48 UTC ·"

$ws.Range("B11").Value = "Yuan Chinês"
$ws.Range("C11").Value = "¥"
$ws.Range("D11").Value = "CNY"
$ws.Range("E11").Value = "0,77"
$ws.Range("F11").Value = "20 de jun., 14:48 UTC ·"

$ws.Range("B12").Value = "Rúpia Indiana"
$ws.Range("C12").Value = "₹"
$ws.Range("D12").Value = "INR"
Set-TextCell $ws "E12" "0,063"
$ws.Range("F12").Value = "20 de jun., 14:49 UTC ·"

$ws.Range("B13").Value = "Peso Chileno"
$ws.Range("C13").Value = "$"
$ws.Range("D13").Value = "CLP"
$ws.Range("E13").Value = "0,0058"
$ws.Range("F13").Value = "20 de jun., 11:19 UTC ·"

$ws.Range("B14").Value = "Peso Argentino"
$ws.Range("C14").Value = "$"
$ws.Range("D14").Value = "ARS"
$ws.Range("E14").Value = "0,0047"
$ws.Range("F14").Value = "20 de jun., 14:48 UTC ·"

$ws.Range("B15").Value = "Peso Colombiano"
$ws.Range("C15").Value = "$"
$ws.Range("D15").Value = "COP"
$ws.Range("E15").Value = "0,0013"
$ws.Range("F15").Value = "20 de jun., 14:48 UTC ·"

$ws.Range("B16").Value = "Rúpia Russa"
$ws.Range("C16").Value = "₽"
$ws.Range("D16").Value = "RUB"
$ws.Range("E16").Value = "1,10"
$ws.Range("F16").Value = "20 de jun., 14:49 UTC ·"

$ws.Range("B17").Value = "Riyal Saudi"
$ws.Range("C17").Value = "﷼"
$ws.Range("D17").Value = "SAR"
$ws.Range("E17").Value = "1,46"
$ws.Range("F17").Value = "20 de jun., 14:49 UTC ·"

$ws.Range("B18").Value = "Dólar de Singapura"
$ws.Range("C18").Value = "S$"
$ws.Range("D18").Value = "SGD"
$ws.Range("E18").Value = "4,27"
$ws.Range("F18").Value = "20 de jun., 14:48 UTC ·"

$ws.Range("B19").Value = "Peso Filipino"
$ws.Range("C19").Value = "₱"
$ws.Range("D19").Value = "PHP"
Set-TextCell $ws "E19" "0,096"
$ws.Range("F19").Value = "20 de jun., 14:49 UTC ·"

$ws.Range("B20").Value = "Yuan de Taiwan"
$ws.Range("C20").Value = "NT$"
$ws.Range("D20").Value = "TWD"
$ws.Range("E20").Value = "4,11"
$ws.Range("F20").Value = "20 de jun., 14:50 UTC ·"

$ws.Range("B21").Value = "Dinar Iraquiano"
$ws.Range("C21").Value = "؋"
$ws.Range("D21").Value = "IQD"
$ws.Range("E21").Value = "0,0042"
$ws.Range("F21").Value = "20 de jun., 14:49 UTC ·"

$ws.Range("B22").Value = "Rúpia Sri Lanka"
$ws.Range("C22").Value = "Rs"
$ws.Range("D22").Value = "LKR"
$ws.Range("E22").Value = "3,47"
$ws.Range("F22").Value = "20 de jun., 14:48 UTC ·"

$ws.Range("B23").Value = "Yuan Chinês"
$ws.Range("C23").Value = "¥"
$ws.Range("D23").Value = "CNY"
$ws.Range("E23").Value = "0,77"
$ws.Range("F23").Value = "20 de jun., 14:48 UTC ·"

$ws.Range("B24").Value = "Won Sul-Coreano"
$ws.Range("C24").Value = "₩"
$ws.Range("D24").Value = "KRW"
$ws.Range("E24").Value = "0,0040"
$ws.Range("F24").Value = "20 de jun., 14:50 UTC ·"
